$wb = $excel.ActiveWorkbook

$mainWs = $wb.Worksheets.Item("Main")
$privateWs = $wb.Worksheets.Item("Private")

# --- Private sheet: new company row plus Raise/Amount columns ---
$privateWs.Range("B11").Value = "Color"
$privateWs.Range("E2").Value = "Raise"
$privateWs.Range("E11").Value = "Series E"
$privateWs.Range("F2").Value = "Amount"
$privateWs.Range("F11").Value = 482

# --- Main sheet: add a couple of new rows below the header ---
$mainWs.Range("A3").Value = "x"
$mainWs.Range("B3").Value = "Oracle"
$mainWs.Range("B4").Value = "Epic"

# Update selection to reflect where the user ended up on each sheet
$privateWs.Range("F11").Select()
$mainWs.Range("B5").Select()

$mainWs.Activate()
